$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.051.15"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.682.33"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "215.79"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "21.43"
$ws.Range("E8").Value = "  +5.22%  "
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.920.36"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.685.54"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "0.534"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "66.10"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.063.67"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "236.74"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").Value = "146.82"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "1.516.63"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").Value = "3.20"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("E35").Value = "  +3.95%  "
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D38").Value = "0.920"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("E40").Value = "  +7.45%  "
$ws.Range("D41").Value = "5.75"
$ws.Range("E41").Value = "  -5.17%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "1.825.62"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "90.22"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "0.105"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").Value = "1.52"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "7.85"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("E51").Value = "  -0.03%  "
